$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Path (B2), Capital (F2), Optimization (G2), and TimeStamp (J2)
# cells to reflect the new uploaded file's details.
$ws.Range("B2").Value = "./NewOnes/CRYPTO Lorenzo Reyes lreyes@udesa.edu.ar 3200 MonteSharpe 2022-11-05.xlsx"

# Capital must stay stored as text (not a number), so force the cell to
# Text format before writing the value.
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "3200"

$ws.Range("G2").Value = "MonteSharpe"
$ws.Range("J2").Value = "21:37:41 05-11-2022"
